$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "How many scales can I define?"
$ws.Range("B29").Value = "llama3.2:latest"
$ws.Range("C29").Value = "According to the document, you can define 23 scales."

$ws.Range("A30").Value = "What the maximum number of data files I can load?"
$ws.Range("B30").Value = "llama3.2:latest"
$ws.Range("C30").Value = "The maximum number of data files you can load is unlimited."
